$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 89.92308
$ws.Range("I6").Value = 20.428572
$ws.Range("J6").Value = 171
$ws.Range("K6").Value = 61.28571599999999
$ws.Range("L6").Value = 513
$ws.Range("M6").Value = 50.71428400000001
$ws.Range("N6").Value = -737

$ws.Range("H107").Value = 2241.9167
$ws.Range("I107").Value = 1499.75
$ws.Range("J107").Value = 3726.25
$ws.Range("K107").Value = 1499.75
$ws.Range("L107").Value = 3726.25
$ws.Range("M107").Value = 420.25
$ws.Range("N107").Value = -7566.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 464
$ws.Range("J2").Value = 532.6667
$ws.Range("L2").Value = 532.6667
$ws.Range("N2").Value = -758.6667

$ws.Range("H32").Value = 2045.3334
$ws.Range("I32").Value = 2057.8
$ws.Range("K32").Value = 2057.8
$ws.Range("M32").Value = -1770.8

$ws.Range("H110").Value = 969.9
$ws.Range("I110").Value = 969.9
$ws.Range("K110").Value = 969.9
$ws.Range("M110").Value = 1075.1

$ws.Range("H116").Value = 464
$ws.Range("J116").Value = 532.6667
$ws.Range("L116").Value = 532.6667
$ws.Range("N116").Value = -5120.6667

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 464
$ws.Range("J3").Value = 532.6667
$ws.Range("L3").Value = 532.6667
$ws.Range("N3").Value = -760.6667

$ws.Range("H20").Value = 3184.5
$ws.Range("I20").Value = 3184.5
$ws.Range("K20").Value = 3184.5
$ws.Range("M20").Value = -2937.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H3").Value = 3990
$ws.Range("I3").Value = 783.3333
$ws.Range("J3").Value = 8800
$ws.Range("K3").Value = 783.3333
$ws.Range("L3").Value = 8800
$ws.Range("M3").Value = -670.3333
$ws.Range("N3").Value = -9026

$ws.Range("H12").Value = 342.14285
$ws.Range("I12").Value = 341
$ws.Range("J12").Value = 345
$ws.Range("K12").Value = 341
$ws.Range("L12").Value = 345
$ws.Range("M12").Value = -171
$ws.Range("N12").Value = -685

$ws.Range("H105").Value = 8799.799999999999
$ws.Range("I105").Value = 1999.5
$ws.Range("J105").Value = 13333.333
$ws.Range("K105").Value = 1999.5
$ws.Range("L105").Value = 13333.333
$ws.Range("M105").Value = -252.5
$ws.Range("N105").Value = -16827.333

$ws.Range("H107").Value = 51249.5
$ws.Range("J107").Value = 51249.5
$ws.Range("L107").Value = 51249.5
$ws.Range("N107").Value = -55089.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 694.875
$ws.Range("I14").Value = 694.875
$ws.Range("K14").Value = 2084.625
$ws.Range("M14").Value = -1911.625

$ws.Range("H68").Value = 1912.6666
$ws.Range("J68").Value = 2115.2
$ws.Range("L68").Value = 6345.599999999999
$ws.Range("N68").Value = -7967.599999999999

$ws.Range("H71").Value = 1912.6666
$ws.Range("J71").Value = 2115.2
$ws.Range("L71").Value = 19036.8
$ws.Range("N71").Value = -27148.8

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 500812.5
$ws.Range("I3").Value = 417500
$ws.Range("K3").Value = 417500
$ws.Range("M3").Value = -417384

$ws.Range("H7").Value = 1255997.2
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 1255997.2
$ws.Range("K7").Value = 0
$ws.Range("L7").Value = 1255997.2
$ws.Range("N7").Value = -1256221.2
$ws.Range("M7").ClearContents()

$ws.Range("H8").Value = 1255997.2
$ws.Range("I8").Value = 0
$ws.Range("J8").Value = 1255997.2
$ws.Range("K8").Value = 0
$ws.Range("L8").Value = 1255997.2
$ws.Range("N8").Value = -1256275.2
$ws.Range("M8").ClearContents()

$ws.Range("H11").Value = 9300000
$ws.Range("I11").Value = 9990909
$ws.Range("J11").Value = 5500000
$ws.Range("K11").Value = 9990909
$ws.Range("L11").Value = 5500000
$ws.Range("M11").Value = -9990770
$ws.Range("N11").Value = -5500278

$ws.Range("H13").Value = 0
$ws.Range("I13").Value = 0
$ws.Range("J13").Value = 0
$ws.Range("K13").Value = 0
$ws.Range("L13").Value = 0
$ws.Range("M13").ClearContents()
$ws.Range("N13").ClearContents()

$ws.Range("H113").Value = 1500
$ws.Range("I113").Value = 1500
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 1500
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 670
$ws.Range("N113").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H3").Value = 11399.75
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 11399.75
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = 11399.75
$ws.Range("N3").Value = -11623.75

$ws.Range("H15").Value = 11399.75
$ws.Range("I15").Value = 0
$ws.Range("J15").Value = 11399.75
$ws.Range("K15").Value = 0
$ws.Range("L15").Value = 11399.75
$ws.Range("N15").Value = -11739.75

$ws.Range("H61").Value = 3367.8333
$ws.Range("I61").Value = 3367.8333
$ws.Range("K61").Value = 3367.8333
$ws.Range("M61").Value = -3165.8333

$ws.Range("H113").Value = 3367.8333
$ws.Range("I113").Value = 3367.8333
$ws.Range("K113").Value = 3367.8333
$ws.Range("M113").Value = -1197.8333

$ws.Range("H136").Value = 8000
$ws.Range("I136").Value = 8000
$ws.Range("K136").Value = 24000
$ws.Range("M136").Value = -21450

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H8").Value = 6009.6665
$ws.Range("I8").Value = 3003
$ws.Range("J8").Value = 7513
$ws.Range("K8").Value = 3003
$ws.Range("L8").Value = 7513
$ws.Range("M8").Value = -2863
$ws.Range("N8").Value = -7793

$ws.Range("H11").Value = 6201.3335
$ws.Range("I11").Value = 1004
$ws.Range("K11").Value = 1004
$ws.Range("M11").Value = -862

$ws.Range("H107").Value = 2199.5
$ws.Range("I107").Value = 1285.1428
$ws.Range("J107").Value = 4333
$ws.Range("K107").Value = 3855.4284
$ws.Range("L107").Value = 12999
$ws.Range("M107").Value = -1935.4284
$ws.Range("N107").Value = -16839

$ws.Range("H113").Value = 987.55554
$ws.Range("I113").Value = 1041.2858
$ws.Range("K113").Value = 3123.8574
$ws.Range("M113").Value = -953.8574000000003
